$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Preguntas Texto Corto"

# Column width for column A
$ws.Columns("A").ColumnWidth = 44.43

# Header row (row 1): question + option headers first (controls the
# shared-string table order)
$ws.Range("A1").Value = "Pregunta"
$ws.Range("B1:O1").Value = "Opcion"

# Header fill: apply the theme color to a single cell, then propagate the
# resulting cell format to the rest of the header row (avoids generating a
# redundant intermediate style for every cell in the range)
$ws.Range("A1").Interior.ThemeColor = 9
$ws.Range("A1").Interior.TintAndShade = 0.79998168889431442
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column A - remaining questions
$ws.Range("A2").Value = "Pregunta 1"
$ws.Range("A3").Value = "Pregunta 2"
$ws.Range("A4").Value = "Pregunta 3"
$ws.Range("A5").Value = "Pregunta 4"
$ws.Range("A6").Value = "Pregunta 5"

# Row 2
$ws.Range("B2").Value = "Opcion 1"
$ws.Range("C2").Value = "Opcion 2"
$ws.Range("D2").Value = "Opcion 3"

# Row 3
$ws.Range("B3").Value = "Opcion 1"
$ws.Range("C3").Value = "Opcion 2"

# Row 4
$ws.Range("B4").Value = "Opcion 1"
$ws.Range("C4").Value = "Opcion 2"
$ws.Range("D4").Value = "Opcion 3"
$ws.Range("E4").Value = "Opcion 4"

# Row 5
$ws.Range("B5").Value = "Opcion 1"

# Row 6
$ws.Range("B6").Value = "Opcion 1"
$ws.Range("C6").Value = "Opcion 2"
$ws.Range("D6").Value = "Opcion 3"
$ws.Range("E6").Value = "Opcion 4"
$ws.Range("F6").Value = "Opcion 5"
$ws.Range("G6").Value = "Opcion 6"

# Selection as in the saved file
$ws.Range("D10").Select() | Out-Null
